# Insert a new column "Rch_id" at the front of the sheet, pushing the
# existing "name"/"age" columns (and their data) one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:B to B:C by inserting a new column before A.
$ws.Columns.Item(1).Insert()

# New header cell for the inserted column.
$ws.Range("A1").Value = "Rch_id"

# Re-select A2 to match the post-edit selection state.
$ws.Range("A2").Select()
